# CIV-17609 updated wording for casenumber
# The template heading "Claim number" should read "Case number".

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "Claim number",  # FindText
    $true,           # MatchCase
    $false,          # MatchWholeWord
    $false,          # MatchWildcards
    $false,          # MatchSoundsLike
    $false,          # MatchAllWordForms
    $true,           # Forward
    1,               # Wrap (wdFindContinue)
    $false,          # Format
    "Case number",   # ReplaceWith
    2                # Replace (wdReplaceAll)
)
